$d = $word.ActiveDocument

# Locate the last bullet item in the symbol list ("[Endor]") so the new
# entries can be appended immediately after it, inheriting the same
# bullet-list paragraph formatting (numId=1, ilvl=0, hanging indent).
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("[Endor]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph '[Endor]'"
}

$hit = $find.Parent
$hitStart = $hit.Start

# Resolve the paragraph index of the found "[Endor]" run so we can keep
# re-fetching live paragraphs by index (Range objects returned from
# InsertParagraphAfter do not auto-track subsequent edits).
$endorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $hitStart) {
        $endorIndex = $i
        break
    }
}

if ($endorIndex -eq 0) {
    throw "Could not resolve paragraph index for anchor '[Endor]'"
}

$items = @("[Megalon]", "[Antioch]", "[Holon]", "[Dalphine]", "[Antorus]")

$idx = $endorIndex
foreach ($item in $items) {
    $d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
    $idx = $idx + 1
    $d.Paragraphs($idx).Range.InsertBefore($item)
}
